# "Installing Python?" -> "Installing Python" (drop the trailing "?")
# on slide 1's "TextBox 6" shape, then let the textbox's auto-fit width
# follow the now-shorter text (spAutoFit / wrap="none" textbox).

$p   = $ppt.ActivePresentation
$s   = $p.Slides.Item(1)
$shp = $s.Shapes.Item(2)          # "TextBox 6" ("Installing Python?")

$shp.TextFrame.TextRange.Text = "Installing Python"

# EMU -> points (1 pt = 12700 EMU); new width matches the re-flowed,
# narrower auto-fit textbox extent (cx 4613764 -> 4501553, cy unchanged).
$shp.Width = 4501553 / 12700
